$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = 1..12
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Range("K22:L22").Select()
